# Task 537: add pricing_prepenalty_allowed and pricing_prepenalty_exists fields
# to the "invalid" worksheet (new columns Q and R).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("Q1").Value = "pricing_prepenalty_allowed"
$ws.Range("R1").Value = "pricing_prepenalty_exists"

# --- Data rows (rows 2-11) --------------------------------------------------
$q = @(1, 2, 999, 0, 3, 1, 2, 999, 0, 0)
$r = @(1, 2, 1, 2, 0, 0, 999, 999, 1, 2)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 17).Value = $q[$i]
    $ws.Cells.Item($row, 18).Value = $r[$i]
}

# --- Formatting: match the wrap-text style used by the other data columns --
$ws.Range("Q1:R11").WrapText = $true

# Header row grows taller to fit the new wrapped column headers.
$ws.Rows.Item(1).RowHeight = 51

# --- View state -------------------------------------------------------------
$ws.Range("R12").Select() | Out-Null
